$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from column I (rows 3-14) into column J
$ws.Range("I3:I14").Copy()
$ws.Range("J3:J14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set the new values for column J
$ws.Range("J4").Value = 2022
$ws.Range("J5").Value = 96.4
$ws.Range("J6").Value = 96.4
$ws.Range("J7").Value = 97.9
$ws.Range("J8").Value = 95.3
$ws.Range("J9").Value = 93.8
$ws.Range("J10").Value = 95.5
$ws.Range("J11").Value = 94.4
$ws.Range("J12").Value = 95
$ws.Range("J13").Value = 98.7
$ws.Range("J14").Value = 97.3

# Update the active selection in the sheet view to match the target
$ws.Range("L10").Select()
